# Adds four new questionnaire participants (11-14) to the
# "Questionnaire Results" sheet, matching the rows appended in the
# commit, then restores the active-cell selection to F7 (the state the
# sheet was left in after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questionnaire Results")

# --- Participant 11: Oliver ---------------------------------------
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Oliver"
$ws.Range("C14").Value = "B"
$ws.Range("D14").Value = "C"
$ws.Range("E14").Value = "B"
$ws.Range("I14").Value = "Tipp B"
$ws.Range("F14").Value = "Tipp B"
$ws.Range("G14").Value = "Shopping"

# --- Participant 12: Hong -------------------------------------------
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "Hong"
$ws.Range("C15").Value = "B"
$ws.Range("D15").Value = "B"
$ws.Range("E15").Value = "C"
$ws.Range("I15").Value = "Tipp A"
$ws.Range("F15").Value = "Tipp B"
$ws.Range("G15").Value = "Going to a playhall to play bowling"

# --- Participant 13: Kemal ------------------------------------------
$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "Kemal"
$ws.Range("C16").Value = "B"
$ws.Range("D16").Value = "A"
$ws.Range("E16").Value = "C"
$ws.Range("I16").Value = "Tipp C"
$ws.Range("F16").Value = "Tipp C"
$ws.Range("G16").Value = "Having dinner together"

# --- Participant 14: Simone -----------------------------------------
$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "Simone"
$ws.Range("C17").Value = "A"
$ws.Range("D17").Value = "A"
$ws.Range("E17").Value = "B"
$ws.Range("I17").Value = "Tipp C"
$ws.Range("F17").Value = "Tipp A"
$ws.Range("G17").Value = "Movie theater"

# Move the active selection to F7, matching the saved view state.
[void]$ws.Range("F7").Select()
